$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: I1 = "I0", J1 = "IF"
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold, border, centered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for new columns I (I0) and J (IF), rows 2-16
$iValues = @(6,9,5,5,5,4,9,6,5,8,7,4,9,7,7)
$jValues = @(9,9,5,7,7,6,9,6,5,8,8,4,9,7,7)

for ($r = 2; $r -le 16; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
